# Rename the paired AHB-diff header columns from the generic "_old"/"_new"
# suffixes to the concrete format-version suffixes "_FV2310"/"_FV2404", turn
# the used range A1:U78 into a native Excel Table ("Table1") with an
# AutoFilter, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header row rename: "<Name>_old" -> "<Name>_FV2310", "<Name>_new" -> "<Name>_FV2404"
$headers = @(
    "Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310",
    "Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
    "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn A1:U78 into a native table ("Table1") with an AutoFilter on the header row.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U78"), $null, 1)
$lo.Name = "Table1"

# 3) Freeze the header row (split/freeze at row 2, i.e. top row frozen).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header row renamed, Table1 created over A1:U78, top row frozen."
